# Commit: "updated af2122 and suis1 dependencies with no gi"
#
# The sheet tab names embedded the full NCBI "gi|...|ref|...|" accession
# string. This strips that down to just the RefSeq accession
# (e.g. "gi|384222553|ref|NC_017250.1|" -> "NC_017250.1").
#
# The author also re-opened/re-saved the workbook with the second sheet
# active (tab selected) and the selection on that sheet moved to A37.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "NC_017250.1"
$ws2.Name = "NC_017251.1"

# Make the second sheet the active/selected tab, matching the saved
# workbook view (activeTab="1" / tabSelected on sheet2).
$ws2.Activate()
$ws2.Range("A37").Select()
